$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the Price column as text so numeric-looking strings
# (e.g. "312.57") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.991.74"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "1.828.76"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").Value = "312.57"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").Value = "0.4604"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").Value = "0.3704"
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("D9").Value = "0.07344"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").Value = "0.8731"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.07966"
$ws.Range("E11").Value = "  +5.06%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "19.78"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.799.10"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "6.578"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "5.343"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "91.50"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "0.000008902"
$ws.Range("E18").Value = "  +2.98%  "
$ws.Range("D19").Value = "1.007"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "14.80"
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("D21").Value = "27.325.85"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "5.118"
$ws.Range("E22").Value = "  -1.71%  "
$ws.Range("D23").Value = "10.55"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "2.128.76"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("D25").Value = "153.04"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").Value = "1.837"
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("D27").Value = "18.34"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "2.046"
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("D29").Value = "5.158"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("D30").Value = "115.53"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").Value = "0.08902"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "2.964"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").Value = "0.7324"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "4.428"
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("D36").Value = "2.464"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "0.01949"
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").Value = "2.938"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").Value = "7.156"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").Value = "0.5145"
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "0.8901"
$ws.Range("E43").Value = "  -11.90%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "0.1628"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "8.189"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.4838"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "1.006"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "10.16"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "102.63"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "1.632"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.06203"
$ws.Range("E51").Value = "  -0.90%  "
